$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (preserves things like leading zeros, trailing zeros, and multi-dot
    # numbers such as '60.325.23' instead of converting to a number),
    # then reset the style to Normal so no stray number-format/style
    # index is left applied to the cell.
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-CellText "D2" "60.325.23"
Set-CellText "E2" "  -2.02%  "
Set-CellText "D3" "3.370.86"
Set-CellText "E3" "  -2.23%  "
Set-CellText "D4" "0.999"
Set-CellText "E4" "  -0.04%  "
Set-CellText "D5" "566.59"
Set-CellText "E5" "  -2.33%  "
Set-CellText "D6" "140.29"
Set-CellText "E6" "  -6.34%  "
Set-CellText "E7" "  +0.05%  "
Set-CellText "D8" "3.371.01"
Set-CellText "E8" "  -2.26%  "
Set-CellText "E9" "  -0.64%  "
Set-CellText "E10" "  -3.98%  "
Set-CellText "E11" "  -2.53%  "
Set-CellText "D12" "0.387"
Set-CellText "E12" "  -1.15%  "
Set-CellText "D13" "3.944.16"
Set-CellText "E13" "  -2.28%  "
Set-CellText "E14" "  +1.08%  "
Set-CellText "D15" "28.00"
Set-CellText "E15" "  +0.06%  "
Set-CellText "D16" "3.372.44"
Set-CellText "E16" "  -2.18%  "
Set-CellText "E17" "  -3.60%  "
Set-CellText "D18" "60.430.85"
Set-CellText "E18" "  -2.05%  "
Set-CellText "D19" "6.17"
Set-CellText "E19" "  -1.90%  "
Set-CellText "D20" "13.80"
Set-CellText "E20" "  -4.11%  "
Set-CellText "D21" "9.00"
Set-CellText "E21" "  -5.33%  "
Set-CellText "D22" "385.46"
Set-CellText "E22" "  -1.17%  "
Set-CellText "E23" "  -2.15%  "
Set-CellText "D24" "73.01"
Set-CellText "E25" "  -0.01%  "
Set-CellText "E26" "  -7.92%  "
Set-CellText "D27" "3.519.79"
Set-CellText "E27" "  -1.90%  "
Set-CellText "E28" "  -1.92%  "
Set-CellText "D30" "7.35"
Set-CellText "E30" "  -5.10%  "
Set-CellText "D31" "7.91"
Set-CellText "E31" "  -4.24%  "
Set-CellText "D32" "2.13"
Set-CellText "E32" "  -2.33%  "
Set-CellText "E33" "  -9.42%  "
Set-CellText "E34" "  +0.02%  "
Set-CellText "D35" "23.49"
Set-CellText "E35" "  -2.37%  "
Set-CellText "D36" "3.400.46"
Set-CellText "E36" "  -1.97%  "
Set-CellText "E37" "  -2.49%  "
Set-CellText "D38" "168.20"
Set-CellText "E38" "  +0.79%  "
Set-CellText "D39" "4.92"
Set-CellText "E39" "  -5.72%  "
Set-CellText "E40" "  -4.68%  "
Set-CellText "D41" "0.0769"
Set-CellText "E41" "  -2.54%  "
Set-CellText "D42" "26.98"
Set-CellText "E42" "  +0.01%  "
Set-CellText "E43" "  -0.04%  "
Set-CellText "E44" "  -2.20%  "
Set-CellText "D45" "4.42"
Set-CellText "E45" "  -1.92%  "
Set-CellText "B46" "Stacks"
Set-CellText "C46" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText "D46" "1.69"
Set-CellText "E46" "  -1.51%  "
Set-CellText "B47" "OKB"
Set-CellText "C47" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-CellText "D47" "41.26"
Set-CellText "E47" "  -2.61%  "
Set-CellText "D48" "2.513.52"
Set-CellText "E48" "  -3.31%  "
Set-CellText "E49" "  -4.65%  "
Set-CellText "D50" "23.15"
Set-CellText "E50" "  -0.23%  "
Set-CellText "E51" "  -2.97%  "
